$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the bottom 5 rows (models 20-24), shrinking the data range from A1:I26 to A1:I21
$ws.Range("A22:I26").Delete() | Out-Null

# Replace model names (column A) and metric values (columns B:I) for rows 2-21
$names = @(
    "model_8_4_19",
    "model_8_4_18",
    "model_8_4_17",
    "model_8_4_11",
    "model_8_4_16",
    "model_8_4_15",
    "model_8_4_14",
    "model_8_4_13",
    "model_8_4_12",
    "model_8_4_0",
    "model_8_4_8",
    "model_8_4_9",
    "model_8_4_10",
    "model_8_4_1",
    "model_8_4_7",
    "model_8_4_2",
    "model_8_4_6",
    "model_8_4_5",
    "model_8_4_3",
    "model_8_4_4"
)

$values = @(
    @(-0.4899014663165022, -51.94316196863086, -2.367791320830892, -6.228032466976575, 1.648881912231445, 5.85093879699707, 4.958014488220215, 5.430737495422363),
    @(-0.4735566462321408, -52.04349923084717, -2.273262592278074, -6.148693175915089, 1.63079297542572, 5.862027645111084, 4.818850517272949, 5.371127128601074),
    @(-0.4690856985143794, -51.96848892047998, -2.255833300999441, -6.126770291371651, 1.625844955444336, 5.853737831115723, 4.793191909790039, 5.354655265808105),
    @(-0.4620836614178456, -47.24303352481616, -2.476963289064361, -5.962697014225828, 1.618095636367798, 5.331510543823242, 5.118735790252686, 5.231379985809326),
    @(-0.45660665494674, -54.43401398414598, -1.978444879511962, -6.06299040057817, 1.612034440040588, 6.126211643218994, 4.384823799133301, 5.306734561920166),
    @(-0.4088928814149293, -55.01727256029157, -1.675318894682547, -5.828908028334003, 1.559229135513306, 6.190669536590576, 3.938566446304321, 5.130857944488525),
    @(-0.3860790362244, -55.53538597759902, -1.510167319901049, -5.716957229247324, 1.533981084823608, 6.247928142547607, 3.695432662963867, 5.046744823455811),
    @(-0.3401426637086478, -53.15639450241662, -1.468558141746838, -5.493355647328267, 1.483143091201782, 5.985017776489258, 3.634176254272461, 4.878743171691895),
    @(-0.3219004001007038, -51.17348532286831, -1.538128839443159, -5.403086797093591, 1.46295428276062, 5.765878677368164, 3.736597061157227, 4.810919761657715),
    @(-0.1691971690486149, -31.89535244377033, -0.5726358269006926, -3.011642962317011, 1.293956756591797, 3.635383605957031, 2.315212249755859, 3.014123439788818),
    @(-0.1172806229979164, -37.87124290857954, -1.32099362360307, -4.167030209475448, 1.236500382423401, 4.295800685882568, 3.41693377494812, 3.882216453552246),
    @(-0.1112317770814055, -39.69899636984894, -1.169432067175858, -4.169607566522189, 1.229806184768677, 4.49779224395752, 3.193807125091553, 3.884153127670288),
    @(-0.1002792800655827, -38.93980212276599, -1.212346952390247, -4.15006023200772, 1.217684864997864, 4.413890838623047, 3.256985902786255, 3.869466066360474),
    @(-0.08269650301470621, -31.21145942477129, -0.4767025449327351, -2.869932654726447, 1.198226094245911, 3.55980396270752, 2.173980474472046, 2.907650232315063),
    @(-0.07062786921059439, -36.81396998641836, -1.107534740002639, -3.887877542169093, 1.184869527816772, 4.178957462310791, 3.102682828903198, 3.672476768493652),
    @(-0.01518691343204259, -30.91061898501704, -0.3801776421163234, -2.757501648783427, 1.123512744903564, 3.526556968688965, 2.031877994537354, 2.823175668716431),
    @(-0.008305798133970033, -34.55237187936577, -0.9013150899433269, -3.521614620146393, 1.115897417068481, 3.929020404815674, 2.799088954925537, 3.397287607192993),
    @(0.04006858261814428, -33.242800971722, -0.6754112596496025, -3.211341987571695, 1.062361240386963, 3.784294605255127, 2.466516494750977, 3.164165735244751),
    @(0.05691350013637586, -29.86791268740071, -0.2685679868519844, -2.573396430958914, 1.043718934059143, 3.411323547363281, 1.867567777633667, 2.684849739074707),
    @(0.09269774251541762, -29.79898450769068, -0.5868629050725029, -2.861517740007113, 1.004116177558899, 3.403706550598145, 2.336156845092773, 2.901327610015869)
)

for ($i = 0; $i -lt $names.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value2 = $names[$i]
    for ($j = 0; $j -lt 8; $j++) {
        $ws.Cells.Item($row, $j + 2).Value2 = $values[$i][$j]
    }
}

Write-Host "Final dimension:" $ws.UsedRange.Address()
